# Apply commit: "updated test triplet controls file column header"
# Adds a new column header "reference_sample" in cell B1 of the first
# worksheet ("gene1"), which is the active/selected sheet in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = "reference_sample"

# Mirror the author's final selection/active cell as recorded in the diff.
$ws.Range("E13").Select()
